# Update the "as_of_utc" timestamp column (AA) for rows 2-26
# on both the "Главные" and "Линейные" worksheets.

$wb = $excel.ActiveWorkbook

$oldValue = "2025-11-13 03:04:09"
$newValue = "2025-11-13 04:11:37"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # Column AA = 27
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
